$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.193.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.649.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.652.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.206.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.271.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +6.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.791.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0976"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
